$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.960.51"
$ws.Range("E2").Value = "  -2.04%  "
$ws.Range("D3").Value = "2.491.48"
$ws.Range("E3").Value = "  -3.38%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'301.13"
$ws.Range("E5").Value = "  -0.35%  "
$ws.Range("D6").Value = "'94.64"
$ws.Range("E6").Value = "  -2.41%  "
$ws.Range("D7").Value = "'0.578"
$ws.Range("E7").Value = "  +0.65%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D9").Value = "'0.528"
$ws.Range("E9").Value = "  -4.09%  "
$ws.Range("D10").Value = "'35.97"
$ws.Range("E10").Value = "  -1.21%  "
$ws.Range("D11").Value = "'0.0803"
$ws.Range("E11").Value = "  -0.85%  "
$ws.Range("D12").Value = "'0.112"
$ws.Range("E12").Value = "  -2.25%  "
$ws.Range("D13").Value = "'7.37"
$ws.Range("E13").Value = "  -3.23%  "
$ws.Range("D14").Value = "2.871.28"
$ws.Range("E14").Value = "  -3.60%  "
$ws.Range("D15").Value = "2.492.70"
$ws.Range("E15").Value = "  -3.29%  "
$ws.Range("D16").Value = "'14.89"
$ws.Range("E16").Value = "  +3.66%  "
$ws.Range("D17").Value = "'0.841"
$ws.Range("E17").Value = "  -5.03%  "
$ws.Range("D18").Value = "42.004.78"
$ws.Range("E18").Value = "  -2.01%  "
$ws.Range("D19").Value = "'12.72"
$ws.Range("E19").Value = "  -1.31%  "
$ws.Range("D20").Value = "0.0₃0963"
$ws.Range("E20").Value = "  -2.67%  "
$ws.Range("D21").Value = "'6.35"
$ws.Range("E21").Value = "  -4.43%  "
$ws.Range("D22").Value = "'70.70"
$ws.Range("E22").Value = "  -1.76%  "
$ws.Range("D23").Value = "'247.04"
$ws.Range("E23").Value = "  -2.96%  "
$ws.Range("E24").Value = "  -2.35%  "
$ws.Range("D25").Value = "'1.99"
$ws.Range("E25").Value = "  -6.17%  "
$ws.Range("D26").Value = "'26.57"
$ws.Range("E26").Value = "  -7.13%  "
$ws.Range("D27").Value = "'0.999"
$ws.Range("E27").Value = "  -0.11%  "
$ws.Range("E28").Value = "  +8.59%  "
$ws.Range("D29").Value = "'10.04"
$ws.Range("E29").Value = "  -1.71%  "
$ws.Range("D30").Value = "'37.07"
$ws.Range("E30").Value = "  -5.40%  "
$ws.Range("D31").Value = "'5.85"
$ws.Range("E31").Value = "  -3.21%  "
$ws.Range("D32").Value = "'153.87"
$ws.Range("E32").Value = "  -0.96%  "
$ws.Range("D33").Value = "'3.26"
$ws.Range("E33").Value = "  -2.62%  "
$ws.Range("E34").Value = "  -5.55%  "
$ws.Range("D35").Value = "'0.0774"
$ws.Range("E35").Value = "  -4.99%  "
$ws.Range("E36").Value = "  -6.21%  "
$ws.Range("D37").Value = "'18.21"
$ws.Range("E37").Value = "  -1.21%  "
$ws.Range("E38").Value = "  -0.99%  "
$ws.Range("E39").Value = "  -1.19%  "
$ws.Range("D40").Value = "'23.86"
$ws.Range("E40").Value = "  +2.01%  "
$ws.Range("D41").Value = "'3.81"
$ws.Range("E41").Value = "  -2.01%  "
$ws.Range("D42").Value = "'3.33"
$ws.Range("E42").Value = "  -2.19%  "
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("D44").Value = "2.040.29"
$ws.Range("E44").Value = "  -1.44%  "
$ws.Range("D45").Value = "'0.0296"
$ws.Range("E45").Value = "  -4.54%  "
$ws.Range("D46").Value = "'1.94"
$ws.Range("E46").Value = "  -7.52%  "
$ws.Range("B47").Value = "BitcoinSV"
$ws.Range("C47").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D47").Value = "'82.83"
$ws.Range("E47").Value = "  -2.79%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").Value = "'8.86"
$ws.Range("E48").Value = "  -3.92%  "
$ws.Range("D49").Value = "2.730.60"
$ws.Range("E49").Value = "  -3.63%  "
$ws.Range("D50").Value = "'71.39"
$ws.Range("E50").Value = "  -5.94%  "
$ws.Range("D51").Value = "'0.186"
$ws.Range("E51").Value = "  -2.15%  "
